$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row additions: P1 = 14, Q1 = 15 (same bold/border/centered style as O1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

for ($r = 2; $r -le 25; $r++) {
    # Observed transform per data row: I:1->2, K:2->1, M:1->2, O:2->1
    $ws.Cells.Item($r, 9).Value  = 2   # I column
    $ws.Cells.Item($r, 11).Value = 1   # K column
    $ws.Cells.Item($r, 13).Value = 2   # M column
    $ws.Cells.Item($r, 15).Value = 1   # O column

    # New columns P and Q, both value 2
    $ws.Cells.Item($r, 16).Value = 2   # P column
    $ws.Cells.Item($r, 17).Value = 2   # Q column
}
